$d = $word.ActiveDocument

function Insert-GuiBefore($paraIndex, $marker) {
    $d.TrackRevisions = $true
    $p = $d.Paragraphs($paraIndex).Range
    $t = $p.Text
    $idx = $t.IndexOf($marker)
    $insertPos = $p.Start + $idx
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertBefore("GUI")
    $d.TrackRevisions = $false
    $rev = $d.Revisions(1)
    $rev.Accept()
}

# 1. Remove the stray "_GoBack" bookmark from its current (empty list-paragraph) location.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2. "a2 C:\a2\edu\btp400\w2017\client\RemoteBankClient.java" -> insert "GUI" before "RemoteBankClient"
Insert-GuiBefore 19 "RemoteBankClient"

# 3. First "Execute "java edu.btp400.w2017.client.RemoteBankClient""
Insert-GuiBefore 38 "RemoteBankClient"

# 4. "xecute "javac edu/btp400/w2017/client/RemoteBankClient.java" -> insert "GUI", then
#    re-insert the "_GoBack" bookmark right after "GUI" (before "RemoteBankClient.java").
Insert-GuiBefore 51 "RemoteBankClient"

$p51 = $d.Paragraphs(51).Range
$t51 = $p51.Text
$idx51 = $t51.IndexOf("RemoteBankClient.java")
$bmPos = $p51.Start + $idx51
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5. Second "Execute "java edu.btp400.w2017.client.RemoteBankClient""
Insert-GuiBefore 71 "RemoteBankClient"

Write-Output "Done"
